$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "егор"
$ws.Range("C2").Value = "губин"
$ws.Range("D2").Value = "выфывфы"
$ws.Range("E2").Value = "выфв"
$ws.Range("F2").Value = "выфвыфв"
$ws.Range("G2").Value = "вфывф"

# Row 3 updates
$ws.Range("B3").Value = "ddsa"
$ws.Range("C3").Value = "dsa"
$ws.Range("D3").Value = "fdss"
$ws.Range("E3").Value = "fdsfs"
$ws.Range("F3").Value = "fdsfds"
$ws.Range("G3").Value = "fdsfsf"

# Row 4 updates
$ws.Range("B4").Value = "Артём"
$ws.Range("C4").Value = "Боков"
$ws.Range("D4").Value = "выфв"
$ws.Range("E4").Value = "выфввыфвф"
$ws.Range("F4").Value = "ыфввф"
$ws.Range("G4").Value = "выфвф"
